# "Generate Report for Archive"
# Update the localization status text from "Ready for handoff" to
# "In Translation" on every sheet that shows it (Overview + per-locale
# sheets), then re-fit the Status column(s) to the new (shorter) text.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns are E (zh-cn) and F (de-de) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"
$wsOverview.Columns.Item(5).ColumnWidth = 12.45
$wsOverview.Columns.Item(6).ColumnWidth = 12.45

# --- zh-cn sheet: status column is C ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"
$wsZhCn.Columns.Item(3).ColumnWidth = 12.45

# --- de-de sheet: status column is C ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"
$wsDeDe.Columns.Item(3).ColumnWidth = 12.45
